$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1096.3684
$ws.Range("I98").Value = 885.05554
$ws.Range("K98").Value = 885.05554
$ws.Range("M98").Value = 612.94446

$ws.Range("H100").Value = 1984
$ws.Range("I100").Value = 1984
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1984
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1443
$ws.Range("N100").ClearContents()

$ws.Range("H122").Value = 1096.3684
$ws.Range("I122").Value = 885.05554
$ws.Range("K122").Value = 2655.16662
$ws.Range("M122").Value = -205.16662

$ws.Range("H125").Value = 8993.25
$ws.Range("I125").Value = 354.33334
$ws.Range("K125").Value = 3189.00006
$ws.Range("M125").Value = -729.0000600000003

$ws.Range("H138").Value = 12765.258
$ws.Range("I138").Value = 3899.4285
$ws.Range("J138").Value = 15351.125
$ws.Range("K138").Value = 11698.2855
$ws.Range("L138").Value = 46053.375
$ws.Range("M138").Value = -6558.2855
$ws.Range("N138").Value = -56333.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3664.1667
$ws.Range("I63").Value = 2996.25
$ws.Range("K63").Value = 2996.25
$ws.Range("M63").Value = -2310.25

$ws.Range("H66").Value = 3664.1667
$ws.Range("I66").Value = 2996.25
$ws.Range("K66").Value = 14981.25
$ws.Range("M66").Value = -11549.25

$ws.Range("H122").Value = 8930527
$ws.Range("I122").Value = 2149.6667
$ws.Range("J122").Value = 25001606
$ws.Range("K122").Value = 6449.000100000001
$ws.Range("L122").Value = 75004818
$ws.Range("M122").Value = -3999.000100000001
$ws.Range("N122").Value = -75009718

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 657.5
$ws.Range("I94").Value = 697.5
$ws.Range("J94").Value = 577.5
$ws.Range("K94").Value = 697.5
$ws.Range("L94").Value = 577.5
$ws.Range("M94").Value = -246.5
$ws.Range("N94").Value = -1479.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3790
$ws.Range("I62").Value = 3790
$ws.Range("K62").Value = 3790
$ws.Range("M62").Value = -3166

$ws.Range("H65").Value = 3790
$ws.Range("I65").Value = 3790
$ws.Range("K65").Value = 3790
$ws.Range("M65").Value = -15830

$ws.Range("H107").Value = 1301.3125
$ws.Range("I107").Value = 2139.1667
$ws.Range("J107").Value = 798.6
$ws.Range("K107").Value = 2139.1667
$ws.Range("L107").Value = 798.6
$ws.Range("M107").Value = -219.1667000000002
$ws.Range("N107").Value = -4638.6

$ws.Range("H122").Value = 6447.5
$ws.Range("I122").Value = 3307.3845
$ws.Range("J122").Value = 20054.666
$ws.Range("K122").Value = 9922.1535
$ws.Range("L122").Value = 60163.99800000001
$ws.Range("M122").Value = -7472.1535
$ws.Range("N122").Value = -65063.99800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 4312.467
$ws.Range("J81").Value = 4399
$ws.Range("L81").Value = 13197
$ws.Range("N81").Value = -15443

$ws.Range("H84").Value = 4312.467
$ws.Range("J84").Value = 4399
$ws.Range("L84").Value = 39591
$ws.Range("N84").Value = -50823

$ws.Range("H123").Value = 2553
$ws.Range("J123").Value = 2663.9285
$ws.Range("L123").Value = 7991.7855
$ws.Range("N123").Value = -12891.7855

$ws.Range("H130").Value = 4646.25
$ws.Range("I130").Value = 3742.5
$ws.Range("J130").Value = 5550
$ws.Range("K130").Value = 11227.5
$ws.Range("L130").Value = 16650
$ws.Range("M130").Value = -6207.5
$ws.Range("N130").Value = -26690

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 28350
$ws.Range("J46").Value = 28350
$ws.Range("L46").Value = 28350
$ws.Range("N46").Value = -28662

$ws.Range("H57").Value = 26153.076
$ws.Range("J57").Value = 26153.076
$ws.Range("L57").Value = 26153.076
$ws.Range("N57").Value = -27793.076

$ws.Range("H97").Value = 1583.75
$ws.Range("I97").Value = 1441.3334
$ws.Range("J97").Value = 2011
$ws.Range("K97").Value = 1441.3334
$ws.Range("L97").Value = 2011
$ws.Range("M97").Value = -945.3334
$ws.Range("N97").Value = -3003

$ws.Range("H106").Value = 9690
$ws.Range("J106").Value = 9690
$ws.Range("L106").Value = 9690
$ws.Range("N106").Value = -12214

$ws.Range("H107").Value = 645.3333
$ws.Range("I107").Value = 622.25
$ws.Range("J107").Value = 656.875
$ws.Range("K107").Value = 622.25
$ws.Range("L107").Value = 656.875
$ws.Range("M107").Value = 1297.75
$ws.Range("N107").Value = -4496.875

$ws.Range("H122").Value = 7156.143
$ws.Range("I122").Value = 10055.143
$ws.Range("K122").Value = 30165.429
$ws.Range("M122").Value = -27715.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2116.6667
$ws.Range("J46").Value = 3000
$ws.Range("L46").Value = 3000
$ws.Range("N46").Value = -3376

$ws.Range("H82").Value = 2099.4285
$ws.Range("I82").Value = 2099.4285
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 2099.4285
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -1738.4285
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 2099.4285
$ws.Range("I85").Value = 2099.4285
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 2099.4285
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -851.4285
$ws.Range("N85").ClearContents()

$ws.Range("H93").Value = 1328.5714
$ws.Range("I93").Value = 975
$ws.Range("K93").Value = 975
$ws.Range("M93").Value = 273

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 23644.5
$ws.Range("I74").Value = 10000
$ws.Range("J74").Value = 26373.4
$ws.Range("K74").Value = 10000
$ws.Range("L74").Value = 26373.4
$ws.Range("M74").Value = -9064
$ws.Range("N74").Value = -28245.4

$ws.Range("H77").Value = 23644.5
$ws.Range("I77").Value = 10000
$ws.Range("J77").Value = 26373.4
$ws.Range("K77").Value = 30000
$ws.Range("L77").Value = 79120.20000000001
$ws.Range("M77").Value = -25320
$ws.Range("N77").Value = -88480.20000000001

$ws.Range("H96").Value = 4700.857
$ws.Range("I96").Value = 1632.6666
$ws.Range("J96").Value = 7002
$ws.Range("K96").Value = 1632.6666
$ws.Range("L96").Value = 7002
$ws.Range("M96").Value = -259.6666
$ws.Range("N96").Value = -9748
